$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$summary1 = @'
In this informational e-book, Megan M. Gunderson describes the key plans and actions taken by President John F. Kennedy during his time in office, highlighting his domestic initiatives and foreign policy challenges through clear and factual language. The text elaborates on significant events such as the Bay of Pigs invasion and the Cuban Missile Crisis, illustrating the complexities of the Cold War era.
'@

$glossed1 = @'
poverty (noun): the state of being extremely poor.
civil rights (noun): the rights of citizens to political and social freedom and equality.
Cold War (proper noun): a state of political hostility between countries characterized by threats, propaganda, and other measures short of open warfare.
Soviet Union (proper noun): a former federal socialist republic in eastern Europe and northern Asia.
Communist (proper adjective): relating to or denoting a political ideology based on communal ownership and the absence of social classes.
space (noun): the physical universe beyond the earth's atmosphere.
relationships (noun): the ways in which two or more people or groups regard and behave toward each other.
Communism (proper noun): a political and economic ideology advocating for a classless system in which the means of production are owned communally.
Bay of Pigs (proper noun): a failed invasion of Cuba by Cuban exiles sponsored by the U.S. government in April 1961.
Communist (proper adjective): relating to or denoting a political ideology based on communal ownership and the absence of social classes (repeated entry from above).
China (proper noun): a country in East Asia.
India (proper noun): a country in South Asia.
weapons (noun): tools or devices used for causing harm or damage, especially in warfare.
Communist (proper adjective): relating to or denoting a political ideology based on communal ownership and the absence of social classes (repeated entry from above).
North Vietnam (proper noun): the communist state in the northern part of Vietnam during the Vietnam War.
advisers (noun): individuals who provide advice or guidance, especially in a professional context.
nuclear (adjective): relating to the nucleus of an atom, often used in the context of nuclear weapons or energy.
Cuban Missile Crisis (proper noun): a 13-day confrontation between the United States and the Soviet Union in October 1962 over the installation of nuclear-armed Soviet missiles in Cuba.
showdown (noun): a decisive confrontation or dispute.
'@

$vocab1 = @'
Making Plans (excerpt from John F. Kennedy)  
By Megan M. Gunderson, Abdo Publishing  
John F. Kennedy became President of the United States on January 20, 1961. He had many plans when he took office. These included programs for higher wages and aid for people experiencing *poverty*. He also hoped to improve *civil rights*.  
During this time, America was part of the *Cold War* with the *Soviet Union*. Another plan was to put a person on the moon before the *Communist* Soviet Union did. As a result, Kennedy greatly improved the US *space* program.  
*Relationships* with other countries were another concern for Kennedy. American officials wanted to stop the spread of *Communism*. So, they trained a group of Cubans to overthrow Cuba's *Communist* leader. In 1961, the group invaded Cuba at the *Bay of Pigs*. However, the effort failed.  
Then, in 1962, *Communist* China invaded India. Kennedy sent *weapons* to India's army. Meanwhile, South Vietnam was fighting *Communist* North Vietnam. Kennedy sent thousands of US military *advisers* to help the South Vietnamese.  
In October 1962, Kennedy learned *nuclear* missile bases were being built in Cuba. He wanted to stop Soviet ships from bringing supplies there. So, Kennedy stated that the US Navy would block Cuba's coast.  
The Soviet leader said his country would guard its shipping rights. A *nuclear* war could begin. But Kennedy held his ground. The Soviet ships stayed away from the blocked area. Kennedy won the *showdown*! This event became known as the *Cuban Missile Crisis*.
'@

$mcq1 = @'
1. What was one of John F. Kennedy's main goals regarding space during the Cold War?
   A. To establish a permanent military base on the moon
   B. To put a person on the moon before the Soviet Union did
   C. To send robotic missions to Mars
   D. To develop a space defense system against the Soviets
   [CCSS.RI.5.1; TEKS.ELAR.5.7(C)]
   Key: B
2. How did the United States attempt to address the spread of Communism in Cuba?
   A. By launching a direct military invasion
   B. By training a group of Cubans to overthrow the government
   C. By negotiating peace treaties with Soviet leaders
   D. By providing financial assistance to the Cuban government
   [CCSS.RI.3.1; TEKS.ELAR.3.7(C)]
   Key: B
3. What action did Kennedy take in response to the discovery of nuclear missile bases in Cuba?
   A. He sent diplomats to negotiate with Soviet leaders
   B. He ordered a blockade of Cuba to prevent more supplies
   C. He initiated a nuclear attack on Cuba
   D. He withdrew all military advisors from the region
   [CCSS.RI.4.1; TEKS.ELAR.4.7(C)]
   Key: B
4. Which event is referred to as the "Cuban Missile Crisis"?
   A. The failed invasion of Cuba at the Bay of Pigs
   B. The construction of nuclear missile bases in Cuba
   C. The conflict between North and South Vietnam
   D. The space race to the moon
   [CCSS.RI.3.1; TEKS.ELAR.3.7(C)]
   Key: B
5. What was a significant outcome of the Cold War on the relationships between the USA and other countries?
   A. Improved trade relations with all nations
   B. Increased military involvement in foreign conflicts
   C. Strengthened alliances with Communist countries
   D. Reduced tension with the Soviet Union
   [CCSS.RI.5.1; TEKS.ELAR.5.7(C)]
   Key: B
'@

$summary2 = @'
In this informational e-book, Megan M. Gunderson describes the significant plans and actions of John F. Kennedy during his presidency, utilizing factual language to outline his domestic and foreign policies. The text highlights key events such as the Bay of Pigs invasion and the Cuban Missile Crisis, providing context for Kennedy's leadership during the Cold War era.
'@

$glossed2 = @'
higher wages (noun): increased payments to employees for their work
poverty (noun): the state of being extremely poor
civil rights (noun): the rights of citizens to political and social freedom and equality
Cold War (proper noun): the geopolitical tension between the Soviet Union and the United States after World War II
Soviet Union (proper noun): a former federal socialist state in Eastern Europe and Northern Asia
Communist (proper adjective): relating to a political party or movement that advocates for communism
space program (noun): a national or international program involving the development and application of space exploration technology
relationships (noun): the way in which two or more concepts, objects, or people are connected
Communism (proper noun): a political and economic ideology advocating for a classless system in which the means of production are communally owned
Bay of Pigs (proper noun): a failed invasion of Cuba by a CIA-sponsored paramilitary group in 1961
Communist China (proper noun): the People's Republic of China under the governance of the Communist Party
India (proper noun): a country in South Asia
South Vietnam (proper noun): the southeastern region of Vietnam that existed as a separate state from 1955 until 1975
Communist North Vietnam (proper noun): the northern part of Vietnam governed by a communist regime
advisers (noun): individuals who provide expert advice or guidance
nuclear (adjective): relating to the atomic nucleus, especially regarding nuclear energy or weapons
Cuban Missile Crisis (proper noun): a 13-day confrontation between the United States and the Soviet Union in 1962 over the presence of missile sites in Cuba
'@

$vocab2 = @'
Making Plans (excerpt from John F. Kennedy)  
By Megan M. Gunderson, Abdo Publishing  
John F. Kennedy became President of the United States on January 20, 1961. He had many plans when he took office. These included programs for *higher wages* and aid for people experiencing *poverty*. He also hoped to improve *civil rights*.  
During this time, America was part of the *Cold War* with the *Soviet Union*. Another plan was to put a person on the moon before the *Communist* Soviet Union did. As a result, Kennedy greatly improved the US *space program*.  
*Relationships* with other countries were another concern for Kennedy. American officials wanted to stop the spread of *Communism*. So, they trained a group of Cubans to overthrow Cuba's *Communist* leader. In 1961, the group invaded Cuba at the *Bay of Pigs*. However, the effort failed.  
Then, in 1962, *Communist* China invaded India. Kennedy sent weapons to India's army. Meanwhile, South Vietnam was fighting *Communist* North Vietnam. Kennedy sent thousands of US military *advisers* to help the South Vietnamese.  
In October 1962, Kennedy learned *nuclear* missile bases were being built in Cuba. He wanted to stop Soviet ships from bringing supplies there. So, Kennedy stated that the US Navy would block Cuba's coast.  
The Soviet leader said his country would guard its shipping rights. A *nuclear* war could begin. But Kennedy held his ground. The Soviet ships stayed away from the blocked area. Kennedy won the *showdown*! This event became known as the *Cuban Missile Crisis*.
'@

$mcq2 = @'
1. What was one of John F. Kennedy's primary goals as President during the Cold War?
   A. To end all military conflicts worldwide  
   B. To put a person on the moon before the Soviet Union  
   C. To eliminate poverty in South America  
   D. To expand the Civil Rights Movement in Europe  
   [CCSS.RI.4.1; TEKS.ELAR.4.7(C); BEST.ELA.K12.EE.3.1]  
   Key: B
2. How did the Bay of Pigs invasion relate to Kennedy's efforts during the Cold War?
   A. It showcased the importance of space exploration.  
   B. It aimed to stop the spread of Communism in Cuba.   
   C. It was a successful overthrow of a Communist leader.  
   D. It was a covert mission to aid South Vietnam.  
   [CCSS.RI.3.1; TEKS.ELAR.3.7(C); CCSS.RI.5.1]  
   Key: B
3. What action did Kennedy take in response to the discovery of nuclear missile bases in Cuba?
   A. He launched a military attack on the bases.  
   B. He sent negotiators to the Soviet Union.  
   C. He ordered a naval blockade of Cuba.  
   D. He withdrew American military advisers from Vietnam.  
   [CCSS.RI.4.1; TEKS.ELAR.4.7(C)]  
   Key: C
4. In what way did the Cold War influence U.S. foreign policy towards Communist states during Kennedy's administration?
   A. The U.S. sought to cultivate friendships with Communist leaders.  
   B. The U.S. aimed to contain and limit the spread of Communism.  
   C. The U.S. planned to ignore conflicts in regions affected by Communism.  
   D. The U.S. focused solely on economic aid to non-Communist nations.  
   [CCSS.RI.3.1; BEST.ELA.K12.EE.3.1; TEKS.ELAR.5.7(C)]  
   Key: B
5. What was a direct consequence of the Cuban Missile Crisis highlighted in Kennedy's actions?
   A. Improved relations between the United States and Soviet Union.  
   B. A significant escalation of the arms race.  
   C. A strategy to support Communism globally.  
   D. A demonstration of American resolve in the face of nuclear threats.  
   [CCSS.RI.5.1; TEKS.ELAR.4.7(C); CCSS.RI.4.1]  
   Key: D
'@

$contentId = "ARP_MakingPlans_G3-5"
$tags = "John F. Kennedy; Cold War; Cuban Missile Crisis; civil rights; space program"

$ws.Range("A29").Value = $contentId
$ws.Range("B29").Value = $summary1
$ws.Range("C29").Value = $tags
$ws.Range("D29").Value = $glossed1
$ws.Range("E29").Value = $vocab1
$ws.Range("F29").Value = $mcq1

$ws.Range("A30").Value = $contentId
$ws.Range("B30").Value = $summary1
$ws.Range("C30").Value = $tags
$ws.Range("D30").Value = $glossed1
$ws.Range("E30").Value = $vocab1
$ws.Range("F30").Value = $mcq1

$ws.Range("A31").Value = $contentId
$ws.Range("B31").Value = $summary2
$ws.Range("C31").Value = $tags
$ws.Range("D31").Value = $glossed2
$ws.Range("E31").Value = $vocab2
$ws.Range("F31").Value = $mcq2

$ws.Range("A32").Value = $contentId
$ws.Range("B32").Value = $summary2
$ws.Range("C32").Value = $tags
$ws.Range("D32").Value = $glossed2
$ws.Range("E32").Value = $vocab2
$ws.Range("F32").Value = $mcq2

$ws.Rows.Item(29).AutoFit()
$ws.Rows.Item(30).AutoFit()
$ws.Rows.Item(31).AutoFit()
$ws.Rows.Item(32).AutoFit()
